$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the 2007年 data row); rows below shift up by one.
$ws.Rows.Item(2).Delete()
